$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.034.84"
$ws.Range("E2").Value = "  -3.64%  "
$ws.Range("D3").Value = "2.528.53"
$ws.Range("E3").Value = "  -3.33%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "544.34"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.98"
$ws.Range("E6").Value = "  -3.65%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("D9").Value = "2.528.28"
$ws.Range("E9").Value = "  -3.32%  "
$ws.Range("E10").Value = "  -2.92%  "
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.359"
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").Value = "2.982.79"
$ws.Range("E14").Value = "  -3.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.69"
$ws.Range("E15").Value = "  -3.73%  "
$ws.Range("D16").Value = "60.160.99"
$ws.Range("E16").Value = "  -3.30%  "
$ws.Range("E17").Value = "  -2.02%  "
$ws.Range("D18").Value = "2.511.14"
$ws.Range("E18").Value = "  -4.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.66"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("E20").Value = "  -2.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.70"
$ws.Range("E21").Value = "  -3.50%  "
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.83"
$ws.Range("E23").Value = "  -4.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.75"
$ws.Range("E24").Value = "  -1.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.449"
$ws.Range("E25").Value = "  -9.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  +1.23%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.643.56"
$ws.Range("E27").Value = "  -3.24%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.163"
$ws.Range("E28").Value = "  -2.51%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.95"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.26"
$ws.Range("E30").Value = "  +1.95%  "
$ws.Range("E31").Value = "  -3.10%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.29"
$ws.Range("E32").Value = "  -4.01%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.84"
$ws.Range("E33").Value = "  -2.63%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "159.37"
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.44"
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.80"
$ws.Range("E37").Value = "  -2.05%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.52"
$ws.Range("E38").Value = "  -3.60%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.69"
$ws.Range("E39").Value = "  -1.94%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.05"
$ws.Range("E40").Value = "  -0.56%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "318.13"
$ws.Range("E41").Value = "  -5.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.81"
$ws.Range("E42").Value = "  -1.75%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "36.84"
$ws.Range("E43").Value = "  -2.16%  "
$ws.Range("B44").Value = "SuiNetwork"
$ws.Range("C44").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.838"
$ws.Range("E44").Value = "  -5.81%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.604"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.81"
$ws.Range("E47").Value = "  -1.48%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.12"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0946"
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("B50").Value = "Hedera"
$ws.Range("C50").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0533"
$ws.Range("E50").Value = "  -2.46%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0233"
$ws.Range("E51").Value = "  -1.78%  "
